# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# ---- Update "last updated" timestamp (row 1 title) ----
$ws.Range("A1").Value = "Datos actualizados a 12 de Agosto de 2020 a las 17:11"

# ---- Swap Montserrat / Islas Malvinas rows (213 <-> 214) ----
# Row 213 now shows Montserrat with Montserrat's stats
$ws.Range("A213").Value = "Montserrat"
$ws.Range("D213").Value = 12
$ws.Range("H213").Value = 1

# Row 214 now shows Islas Malvinas with Islas Malvinas's stats
$ws.Range("A214").Value = "Islas Malvinas"
$ws.Range("D214").Value = 13
$ws.Range("H214").Value = 0

# ---- Update country statistics (Covid-19 counts) ----

# Row 4: Estados Unidos
$ws.Range("B4").Value = 5318530
$ws.Range("C4").Value = 12573
$ws.Range("D4").Value = 2757622
$ws.Range("E4").Value = 2392896
$ws.Range("G4").Value = 263
$ws.Range("H4").Value = 168012

# Row 19: Argentina
$ws.Range("D19").Value = 187283
$ws.Range("E19").Value = 68540
$ws.Range("G19").Value = 84
$ws.Range("H19").Value = 5088

# Row 22: Alemania
$ws.Range("B22").Value = 219648
$ws.Range("C22").Value = 118
$ws.Range("E22").Value = 10479
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 9269

# Row 23: Francia
$ws.Range("D23").Value = 83237
$ws.Range("E23").Value = 90581

# Row 64: Moldavia
$ws.Range("B64").Value = 28697
$ws.Range("C64").Value = 474
$ws.Range("D64").Value = 19998
$ws.Range("E64").Value = 7836
$ws.Range("G64").Value = 6
$ws.Range("H64").Value = 863

# Row 93: Tayikistan
$ws.Range("B93").Value = 7912
$ws.Range("C93").Value = 41
$ws.Range("D93").Value = 6696
$ws.Range("E93").Value = 1153

# Row 170: Birmania
$ws.Range("B170").Value = 361
$ws.Range("C170").Value = 1
$ws.Range("D170").Value = 318
$ws.Range("E170").Value = 37

# Row 176: Trinidad y Tobago
$ws.Range("B176").Value = 308
$ws.Range("C176").Value = 8
$ws.Range("E176").Value = 161
